$wb = $excel.ActiveWorkbook

$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Copy the test case rows (A2:B9) from "Tests" into "Result" (A2:B9)
for ($r = 2; $r -le 9; $r++) {
    $wsResult.Cells.Item($r, 1).Value2 = $wsTests.Cells.Item($r, 1).Value2
    $wsResult.Cells.Item($r, 2).Value2 = $wsTests.Cells.Item($r, 2).Value2
}

# Update selections on each sheet
$wsTests.Range("B1").Select()
$wsResult.Range("B5").Select()

# Make "Result" the active/visible tab
$wsResult.Activate()

$wb.Save()
